# Auto-generated PowerShell Excel COM-interop script
# Applies updated TPM-derived values to Psen1-Notch2 LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.981972
$ws.Range("H2").Value = 47.945916
$ws.Range("I2").Value = 0.1372507760882863
$ws.Range("J2").Value = 0.1372507760882863
$ws.Range("M2").Value = 2.231113333333334
$ws.Range("N2").Value = 6.69334
$ws.Range("O2").Value = 0.01598125358798882
$ws.Range("P2").Value = 0.01598125358798882
$ws.Range("Q2").Value = 35.65759082216
$ws.Range("R2").Value = 320.91831739944
$ws.Range("S2").Value = 0.002193439457815175
$ws.Range("T2").Value = 0.002193439457815175
$ws.Range("G3").Value = 15.981972
$ws.Range("H3").Value = 47.945916
$ws.Range("I3").Value = 0.1372507760882863
$ws.Range("J3").Value = 0.1372507760882863
$ws.Range("O3").Value = 0.1634493267640196
$ws.Range("P3").Value = 0.1634493267640195
$ws.Range("Q3").Value = 364.6903656099479
$ws.Range("R3").Value = 3282.213290489532
$ws.Range("S3").Value = 0.02243354694946958
$ws.Range("T3").Value = 0.02243354694946958
$ws.Range("G4").Value = 15.981972
$ws.Range("H4").Value = 47.945916
$ws.Range("I4").Value = 0.1372507760882863
$ws.Range("J4").Value = 0.1372507760882863
$ws.Range("M4").Value = 58.02175166666666
$ws.Range("N4").Value = 174.065255
$ws.Range("O4").Value = 0.4156043142904646
$ws.Range("P4").Value = 0.4156043142904646
$ws.Range("Q4").Value = 927.3020105276198
$ws.Range("R4").Value = 8345.718094748578
$ws.Range("S4").Value = 0.0570420146820063
$ws.Range("T4").Value = 0.0570420146820063
$ws.Range("G5").Value = 15.981972
$ws.Range("H5").Value = 47.945916
$ws.Range("I5").Value = 0.1372507760882863
$ws.Range("J5").Value = 0.1372507760882863
$ws.Range("M5").Value = 15.16934033333333
$ws.Range("N5").Value = 45.508021
$ws.Range("O5").Value = 0.1086565487318021
$ws.Range("P5").Value = 0.1086565487318021
$ws.Range("Q5").Value = 242.435972465804
$ws.Range("R5").Value = 2181.923752192236
$ws.Range("S5").Value = 0.01491319564051454
$ws.Range("T5").Value = 0.01491319564051454
$ws.Range("G6").Value = 15.981972
$ws.Range("H6").Value = 47.945916
$ws.Range("I6").Value = 0.1372507760882863
$ws.Range("J6").Value = 0.1372507760882863
$ws.Range("M6").Value = 41.36709099999999
$ws.Range("N6").Value = 124.101273
$ws.Range("O6").Value = 0.2963085566257249
$ws.Range("P6").Value = 0.2963085566257249
$ws.Range("Q6").Value = 661.1276900834519
$ws.Range("R6").Value = 5950.149210751068
$ws.Range("S6").Value = 0.04066857935848066
$ws.Range("T6").Value = 0.04066857935848066
$ws.Range("H7").Value = 69.213024
$ws.Range("I7").Value = 0.1981303529463737
$ws.Range("J7").Value = 0.1981303529463737
$ws.Range("M7").Value = 2.231113333333334
$ws.Range("N7").Value = 6.69334
$ws.Range("O7").Value = 0.01598125358798882
$ws.Range("P7").Value = 0.01598125358798882
$ws.Range("Q7").Value = 51.47403356224001
$ws.Range("R7").Value = 463.2663020601601
$ws.Range("S7").Value = 0.003166371413913725
$ws.Range("T7").Value = 0.003166371413913725
$ws.Range("H8").Value = 69.213024
$ws.Range("I8").Value = 0.1981303529463737
$ws.Range("J8").Value = 0.1981303529463737
$ws.Range("O8").Value = 0.1634493267640196
$ws.Range("P8").Value = 0.1634493267640195
$ws.Range("Q8").Value = 526.4540785398721
$ws.Range("R8").Value = 4738.086706858848
$ws.Range("S8").Value = 0.03238427280060236
$ws.Range("T8").Value = 0.03238427280060235
$ws.Range("H9").Value = 69.213024
$ws.Range("I9").Value = 0.1981303529463737
$ws.Range("J9").Value = 0.1981303529463737
$ws.Range("M9").Value = 58.02175166666666
$ws.Range("N9").Value = 174.065255
$ws.Range("O9").Value = 0.4156043142904646
$ws.Range("P9").Value = 0.4156043142904646
$ws.Range("Q9").Value = 1338.62029687568
$ws.Range("R9").Value = 12047.58267188112
$ws.Range("S9").Value = 0.08234382947640537
$ws.Range("T9").Value = 0.08234382947640535
$ws.Range("H10").Value = 69.213024
$ws.Range("I10").Value = 0.1981303529463737
$ws.Range("J10").Value = 0.1981303529463737
$ws.Range("M10").Value = 15.16934033333333
$ws.Range("N10").Value = 45.508021
$ws.Range("O10").Value = 0.1086565487318021
$ws.Range("P10").Value = 0.1086565487318021
$ws.Range("Q10").Value = 349.971972185056
$ws.Range("R10").Value = 3149.747749665504
$ws.Range("S10").Value = 0.02152816035016681
$ws.Range("T10").Value = 0.0215281603501668
$ws.Range("H11").Value = 69.213024
$ws.Range("I11").Value = 0.1981303529463737
$ws.Range("J11").Value = 0.1981303529463737
$ws.Range("M11").Value = 41.36709099999999
$ws.Range("N11").Value = 124.101273
$ws.Range("O11").Value = 0.2963085566257249
$ws.Range("P11").Value = 0.2963085566257249
$ws.Range("Q11").Value = 954.380487397728
$ws.Range("R11").Value = 8589.424386579552
$ws.Range("S11").Value = 0.05870771890528543
$ws.Range("T11").Value = 0.05870771890528542
$ws.Range("G12").Value = 40.09539033333333
$ws.Range("H12").Value = 120.286171
$ws.Range("I12").Value = 0.3443331924754199
$ws.Range("J12").Value = 0.3443331924754199
$ws.Range("M12").Value = 2.231113333333334
$ws.Range("N12").Value = 6.69334
$ws.Range("O12").Value = 0.01598125358798882
$ws.Range("P12").Value = 0.01598125358798882
$ws.Range("Q12").Value = 89.45735997790446
$ws.Range("R12").Value = 805.11623980114
$ws.Range("S12").Value = 0.005502876067711449
$ws.Range("T12").Value = 0.005502876067711449
$ws.Range("G13").Value = 40.09539033333333
$ws.Range("H13").Value = 120.286171
$ws.Range("I13").Value = 0.3443331924754199
$ws.Range("J13").Value = 0.3443331924754199
$ws.Range("O13").Value = 0.1634493267640196
$ws.Range("P13").Value = 0.1634493267640195
$ws.Range("Q13").Value = 914.9310585662963
$ws.Range("R13").Value = 8234.379527096666
$ws.Range("S13").Value = 0.05628102849261295
$ws.Range("T13").Value = 0.05628102849261294
$ws.Range("G14").Value = 40.09539033333333
$ws.Range("H14").Value = 120.286171
$ws.Range("I14").Value = 0.3443331924754199
$ws.Range("J14").Value = 0.3443331924754199
$ws.Range("M14").Value = 58.02175166666666
$ws.Range("N14").Value = 174.065255
$ws.Range("O14").Value = 0.4156043142904646
$ws.Range("P14").Value = 0.4156043142904646
$ws.Range("Q14").Value = 2326.404780898734
$ws.Range("R14").Value = 20937.6430280886
$ws.Range("S14").Value = 0.1431063603461934
$ws.Range("T14").Value = 0.1431063603461934
$ws.Range("G15").Value = 40.09539033333333
$ws.Range("H15").Value = 120.286171
$ws.Range("I15").Value = 0.3443331924754199
$ws.Range("J15").Value = 0.3443331924754199
$ws.Range("M15").Value = 15.16934033333333
$ws.Range("N15").Value = 45.508021
$ws.Range("O15").Value = 0.1086565487318021
$ws.Range("P15").Value = 0.1086565487318021
$ws.Range("Q15").Value = 608.2206217641767
$ws.Range("R15").Value = 5473.985595877591
$ws.Range("S15").Value = 0.03741405630818246
$ws.Range("T15").Value = 0.03741405630818246
$ws.Range("G16").Value = 40.09539033333333
$ws.Range("H16").Value = 120.286171
$ws.Range("I16").Value = 0.3443331924754199
$ws.Range("J16").Value = 0.3443331924754199
$ws.Range("M16").Value = 41.36709099999999
$ws.Range("N16").Value = 124.101273
$ws.Range("O16").Value = 0.2963085566257249
$ws.Range("P16").Value = 0.2963085566257249
$ws.Range("Q16").Value = 1658.62966059952
$ws.Range("R16").Value = 14927.66694539568
$ws.Range("S16").Value = 0.1020288712607196
$ws.Range("T16").Value = 0.1020288712607196
$ws.Range("G17").Value = 8.831340666666666
$ws.Range("H17").Value = 26.494022
$ws.Range("I17").Value = 0.07584222775512579
$ws.Range("J17").Value = 0.07584222775512579
$ws.Range("M17").Value = 2.231113333333334
$ws.Range("N17").Value = 6.69334
$ws.Range("O17").Value = 0.01598125358798882
$ws.Range("P17").Value = 0.01598125358798882
$ws.Range("Q17").Value = 19.70372191260889
$ws.Range("R17").Value = 177.33349721348
$ws.Range("S17").Value = 0.001212053874432669
$ws.Range("T17").Value = 0.001212053874432669
$ws.Range("G18").Value = 8.831340666666666
$ws.Range("H18").Value = 26.494022
$ws.Range("I18").Value = 0.07584222775512579
$ws.Range("J18").Value = 0.07584222775512579
$ws.Range("O18").Value = 0.1634493267640196
$ws.Range("P18").Value = 0.1634493267640195
$ws.Range("Q18").Value = 201.5211174536327
$ws.Range("R18").Value = 1813.690057082694
$ws.Range("S18").Value = 0.01239636106685875
$ws.Range("T18").Value = 0.01239636106685875
$ws.Range("G19").Value = 8.831340666666666
$ws.Range("H19").Value = 26.494022
$ws.Range("I19").Value = 0.07584222775512579
$ws.Range("J19").Value = 0.07584222775512579
$ws.Range("M19").Value = 58.02175166666666
$ws.Range("N19").Value = 174.065255
$ws.Range("O19").Value = 0.4156043142904646
$ws.Range("P19").Value = 0.4156043142904646
$ws.Range("Q19").Value = 512.4098550450676
$ws.Range("R19").Value = 4611.688695405609
$ws.Range("S19").Value = 0.0315203570604303
$ws.Range("T19").Value = 0.0315203570604303
$ws.Range("G20").Value = 8.831340666666666
$ws.Range("H20").Value = 26.494022
$ws.Range("I20").Value = 0.07584222775512579
$ws.Range("J20").Value = 0.07584222775512579
$ws.Range("M20").Value = 15.16934033333333
$ws.Range("N20").Value = 45.508021
$ws.Range("O20").Value = 0.1086565487318021
$ws.Range("P20").Value = 0.1086565487318021
$ws.Range("Q20").Value = 133.9656121722735
$ws.Range("R20").Value = 1205.690509550462
$ws.Range("S20").Value = 0.008240754716003261
$ws.Range("T20").Value = 0.008240754716003261
$ws.Range("G21").Value = 8.831340666666666
$ws.Range("H21").Value = 26.494022
$ws.Range("I21").Value = 0.07584222775512579
$ws.Range("J21").Value = 0.07584222775512579
$ws.Range("M21").Value = 41.36709099999999
$ws.Range("N21").Value = 124.101273
$ws.Range("O21").Value = 0.2963085566257249
$ws.Range("P21").Value = 0.2963085566257249
$ws.Range("Q21").Value = 365.3268730100006
$ws.Range("R21").Value = 3287.941857090006
$ws.Range("S21").Value = 0.02247270103740082
$ws.Range("T21").Value = 0.02247270103740082
$ws.Range("G22").Value = 28.463871
$ws.Range("H22").Value = 85.39161300000001
$ws.Range("I22").Value = 0.2444434507347945
$ws.Range("J22").Value = 0.2444434507347945
$ws.Range("M22").Value = 2.231113333333334
$ws.Range("N22").Value = 6.69334
$ws.Range("O22").Value = 0.01598125358798882
$ws.Range("P22").Value = 0.01598125358798882
$ws.Range("Q22").Value = 63.50612210638
$ws.Range("R22").Value = 571.55509895742
$ws.Range("S22").Value = 0.003906512774115802
$ws.Range("T22").Value = 0.003906512774115803
$ws.Range("G23").Value = 28.463871
$ws.Range("H23").Value = 85.39161300000001
$ws.Range("I23").Value = 0.2444434507347945
$ws.Range("J23").Value = 0.2444434507347945
$ws.Range("O23").Value = 0.1634493267640196
$ws.Range("P23").Value = 0.1634493267640195
$ws.Range("Q23").Value = 649.513058943189
$ws.Range("R23").Value = 5845.617530488701
$ws.Range("S23").Value = 0.03995411745447595
$ws.Range("T23").Value = 0.03995411745447594
$ws.Range("G24").Value = 28.463871
$ws.Range("H24").Value = 85.39161300000001
$ws.Range("I24").Value = 0.2444434507347945
$ws.Range("J24").Value = 0.2444434507347945
$ws.Range("M24").Value = 58.02175166666666
$ws.Range("N24").Value = 174.065255
$ws.Range("O24").Value = 0.4156043142904646
$ws.Range("P24").Value = 0.4156043142904646
$ws.Range("Q24").Value = 1651.523654634035
$ws.Range("R24").Value = 14863.71289170631
$ws.Range("S24").Value = 0.1015917527254292
$ws.Range("T24").Value = 0.1015917527254292
$ws.Range("G25").Value = 28.463871
$ws.Range("H25").Value = 85.39161300000001
$ws.Range("I25").Value = 0.2444434507347945
$ws.Range("J25").Value = 0.2444434507347945
$ws.Range("M25").Value = 15.16934033333333
$ws.Range("N25").Value = 45.508021
$ws.Range("O25").Value = 0.1086565487318021
$ws.Range("P25").Value = 0.1086565487318021
$ws.Range("Q25").Value = 431.778146403097
$ws.Range("R25").Value = 3886.003317627873
$ws.Range("S25").Value = 0.02656038171693507
$ws.Range("T25").Value = 0.02656038171693507
$ws.Range("G26").Value = 28.463871
$ws.Range("H26").Value = 85.39161300000001
$ws.Range("I26").Value = 0.2444434507347945
$ws.Range("J26").Value = 0.2444434507347945
$ws.Range("M26").Value = 41.36709099999999
$ws.Range("N26").Value = 124.101273
$ws.Range("O26").Value = 0.2963085566257249
$ws.Range("P26").Value = 0.2963085566257249
$ws.Range("Q26").Value = 1177.467541869261
$ws.Range("R26").Value = 10597.20787682335
$ws.Range("S26").Value = 0.07243068606383846
$ws.Range("T26").Value = 0.07243068606383846
